$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 40, pushing the existing rows 40-42 down to 41-43.
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new weekly data point.
$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(40, 3).Value = "Coquimbo"
$ws.Cells.Item(40, 4).Value = 45194
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(40, 6).Value = 100112013
$ws.Cells.Item(40, 7).Value = "Alcachofa"
$ws.Cells.Item(40, 8).Value = "Española"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 600
$ws.Cells.Item(40, 11).Value = 8500
$ws.Cells.Item(40, 12).Value = 9000
$ws.Cells.Item(40, 13).Value = 8750
$ws.Cells.Item(40, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 292
$ws.Cells.Item(40, 17).Value = 30
$ws.Cells.Item(40, 18).Value = "Hortaliza"
